$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Лист1")

# Fill in the "Chapter 8" pomodoro tracking row (row 8) with the missing
# measurements. The I8 (H8-G8) and K8 (I8/J2) formulas, together with all
# the dependent MEDIAN/AVERAGE/SUM summary cells in column M, recalculate
# automatically once the shared formulas in I2:I7 / K2:K7 pick up row 8.
$ws.Range("G8").Value = 221
$ws.Range("H8").Value = 300
$ws.Range("I8").Formula = "=H8-G8"
$ws.Range("J8").Value = 12
$ws.Range("K8").Formula = "=I8/J8"

# Restore the selected cell as recorded in the saved workbook view.
$ws.Range("K13").Select() | Out-Null

$excel.Calculate() | Out-Null
